# Updates the crypto price/volume table (columns D and E, rows 2-51) with
# freshly scraped values, mirroring the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "29.228.70"; E = "-3.86%" },
    @{ Row = 3; D = "1.964.98"; E = "-6.70%" },
    @{ Row = 4; D = "1.015"; E = "+1.38%" },
    @{ Row = 5; D = "328.43"; E = "-4.37%" },
    @{ Row = 6; D = "1.015"; E = "+1.41%" },
    @{ Row = 7; D = "0.4986"; E = "-6.31%" },
    @{ Row = 8; D = "0.4209"; E = "-5.25%" },
    @{ Row = 9; D = "53.65"; E = "-2.22%" },
    @{ Row = 10; D = "0.08906"; E = "-5.38%" },
    @{ Row = 11; D = "1.100"; E = "-6.18%" },
    @{ Row = 12; D = "23.02"; E = "-7.41%" },
    @{ Row = 13; D = "1.969.69"; E = "-4.91%" },
    @{ Row = 14; D = "7.878"; E = "-8.09%" },
    @{ Row = 15; D = "6.410"; E = "-7.46%" },
    @{ Row = 16; D = "1.017"; E = "+1.48%" },
    @{ Row = 17; D = "0.00001103"; E = "-5.10%" },
    @{ Row = 18; D = "91.72"; E = "-9.99%" },
    @{ Row = 19; D = "0.06714"; E = "+0.40%" },
    @{ Row = 20; D = "19.29"; E = "-9.34%" },
    @{ Row = 21; D = "1.016"; E = "+1.55%" },
    @{ Row = 22; D = "5.929"; E = "-6.39%" },
    @{ Row = 23; D = "29.257.71"; E = "-3.82%" },
    @{ Row = 24; D = "11.90"; E = "-5.25%" },
    @{ Row = 25; D = "2.308"; E = "-0.36%" },
    @{ Row = 26; D = "20.65"; E = "-5.72%" },
    @{ Row = 27; D = "155.43"; E = "-4.52%" },
    @{ Row = 28; D = "6.213"; E = "-8.95%" },
    @{ Row = 29; D = "2.296"; E = "-9.16%" },
    @{ Row = 30; D = "127.01"; E = "-5.15%" },
    @{ Row = 31; D = "1.052"; E = "-8.41%" },
    @{ Row = 32; D = "0.09879"; E = "-6.47%" },
    @{ Row = 33; D = "1.515"; E = "-9.48%" },
    @{ Row = 34; D = "5.796"; E = "-7.55%" },
    @{ Row = 35; D = "3.745"; E = "-2.70%" },
    @{ Row = 36; D = "0.02428"; E = "-8.50%" },
    @{ Row = 37; D = "9.150"; E = "-10.43%" },
    @{ Row = 38; D = "0.06334"; E = "-6.97%" },
    @{ Row = 39; D = "1.287"; E = "-4.54%" },
    @{ Row = 40; D = "0.6478"; E = "-7.85%" },
    @{ Row = 41; D = ""; E = "-9.39%" },
    @{ Row = 42; D = "0.2021"; E = "-9.23%" },
    @{ Row = 43; D = "1.015"; E = "+1.51%" },
    @{ Row = 44; D = "0.6261"; E = "-8.97%" },
    @{ Row = 45; D = "13.43"; E = "-7.48%" },
    @{ Row = 46; D = "2.187"; E = "-6.89%" },
    @{ Row = 47; D = "1.281"; E = "-7.16%" },
    @{ Row = 48; D = "3.486"; E = "-4.03%" },
    @{ Row = 49; D = "0.00000000339"; E = "-3.05%" },
    @{ Row = 50; D = "0.06863"; E = "-5.58%" },
    @{ Row = 51; D = "1.117"; E = "-8.73%" }
)

foreach ($u in $updates) {
    if ($u.D -ne "") {
        $dCell = $ws.Range("D" + $u.Row)
        # Price strings such as "1.015" or "0.4986" are valid numeric
        # literals; force the cell to Text first so Excel stores them as
        # strings (matching the source data) instead of coercing to numbers.
        # Values like "29.228.70" (multiple dots) are never numeric and are
        # left with the default formatting, same as the original file.
        if ($u.D -match '^[+-]?(\d+\.)?\d+$') {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }
    $ws.Range("E" + $u.Row).Value = "  " + $u.E + "  "
}
